$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A for "Line #" (everything shifts right by one)
$ws.Columns.Item(1).Insert()

# Insert a new column at D for "Designator" (old D "ModuleAssembly" -> E, etc.)
$ws.Columns.Item(4).Insert()

$ws.Range("A1").Value = "Column=Line #"
$ws.Range("D1").Value = "Column=Designator"

# A1 has no left-hand neighbor for Excel to copy formatting from automatically,
# so explicitly copy the header format (fill/border/alignment) from B1.
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D2").Select()
